$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# UML Design / Implementation subtasks for the "Role Creation" and
# "Role Management" iterations are now complete -- zero out remaining work.
$ws.Range("F14").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("F18").Value = 0
$ws.Range("F19").Value = 0

# Move the active selection to reflect where the user finished editing.
$ws.Range("F18").Select()
